$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: rename to short English field names ---
$ws.Range("A1").Value = "mx_state"
$ws.Range("B1").Value = "mx_municipality"
$ws.Range("C1").Value = "n_matriculas"
$ws.Range("D1").Value = "pct_matriculas"

# --- Title-case Spanish connector words (de/del/la/el/los/y) in place names ---
$titleCaseFixes = @{
    "B5" = "Pabellón De Arteaga"
    "B6" = "Rincón De Romos"
    "B7" = "San Francisco De Los Romo"
    "B25" = "Comitán De Domínguez"
    "B33" = "Marqués De Comillas"
    "B34" = "Mazapa De Madero"
    "B36" = "Montecristo De Guerrero"
    "B63" = "Guadalupe Y Calvo"
    "B64" = "Hidalgo Del Parral"
    "B78" = "San Juan De Sabinas"
    "A88" = "Ciudad De México"
    "B105" = "Coneto De Comonfort"
    "B117" = "Nombre De Dios"
    "B120" = "San Juan De Guadalupe"
    "A127" = "Estado De México"
    "B127" = "Acambay De Ruíz Castañeda"
    "B129" = "Almoloya De Juárez"
    "B137" = "Ecatepec De Morelos"
    "B139" = "Ixtapan De La Sal"
    "B145" = "Naucalpan De Juárez"
    "B147" = "San Felipe Del Progreso"
    "B154" = "Tlalnepantla De Baz"
    "B166" = "Apaseo El Alto"
    "B172" = "Dolores Hidalgo Cuna De La Independencia Nacional"
    "B180" = "Purísima Del Rincón"
    "B185" = "San Francisco Del Rincón"
    "B186" = "San Luis De La Paz"
    "B187" = "Santa Cruz De Juventino Rosas"
    "B188" = "Silao De La Victoria"
    "B192" = "Valle De Santiago"
    "B195" = "Acapulco De Juárez"
    "B196" = "Alcozauca De Guerrero"
    "B199" = "Atoyac De Álvarez"
    "B200" = "Ayutla De Los Libres"
    "B203" = "Chilapa De Álvarez"
    "B204" = "Chilpancingo De Los Bravo"
    "B206" = "Coyuca De Benítez"
    "B207" = "Coyuca De Catalán"
    "B210" = "Cuetzala Del Progreso"
    "B215" = "Iguala De La Independencia"
    "B217" = "Ixcateopan De Cuauhtémoc"
    "B218" = "Zihuatanejo De Azueta"
    "B220" = "La Unión De Isidoro Montes De Oca"
    "B223" = "Mártir De Cuilapan"
    "B230" = "Taxco De Alarcón"
    "B232" = "Técpan De Galeana"
    "B234" = "Tepecoacuilco De Trujano"
    "B237" = "Tlalixtaquilla De Maldonado"
    "B238" = "Tlapa De Comonfort"
    "B245" = "Agua Blanca De Iturbide"
    "B249" = "Cuautepec De Hinojosa"
    "B250" = "Huasca De Ocampo"
    "B254" = "Omitlán De Juárez"
    "B255" = "Pachuca De Soto"
    "B257" = "Progreso De Obregón"
    "B260" = "Santiago De Anaya"
    "B261" = "Santiago Tulantepec De Lugo Guerrero"
    "B263" = "Tenango De Doria"
    "B266" = "Tula De Allende"
    "B267" = "Tulancingo De Bravo"
    "B268" = "Zacualtipán De Ángeles"
    "B277" = "Autlán De Navarro"
    "B290" = "Huejuquilla El Alto"
    "B291" = "Ixtlahuacán Del Río"
    "B294" = "Jilotlán De Los Dolores"
    "B298" = "Lagos De Moreno"
    "B302" = "Ojuelos De Jalisco"
    "B307" = "San Diego De Alejandría"
    "B309" = "San Miguel El Alto"
    "B310" = "San Sebastián Del Oeste"
    "B314" = "Tepatitlán De Morelos"
    "B316" = "Tizapán El Alto"
    "B317" = "Tlajomulco De Zúñiga"
    "B323" = "Unión De San Antonio"
    "B324" = "Unión De Tula"
    "B325" = "Valle De Juárez"
    "B327" = "Yahualica De González Gallo"
    "B342" = "Cojumatlán De Régules"
    "B411" = "Tlaltizapán De Zapata"
    "B418" = "Amatlán De Cañas"
    "B419" = "Bahía De Banderas"
    "B432" = "Acatlán De Pérez Figueroa"
    "B434" = "Heroica Ciudad De Ejutla De Crespo"
    "B435" = "Heroica Ciudad De Huajuapan De León"
    "B440" = "Miahuatlán De Porfirio Díaz"
    "B442" = "Oaxaca De Juárez"
    "B443" = "Ocotlán De Morelos"
    "B444" = "Putla Villa De Guerrero"
    "B447" = "San Antonino El Alto"
    "B462" = "San Miguel Del Puerto"
    "B463" = "San Miguel El Grande"
    "B493" = "Tamazulápam Del Espíritu Santo"
    "B494" = "Tataltepec De Valdés"
    "B495" = "Villa De Chilapa De Díaz"
    "B496" = "Villa De Tututepec De Melchor Ocampo"
    "B512" = "Izúcar De Matamoros"
    "B513" = "Los Reyes De Juárez"
    "B515" = "Palmar De Bravo"
    "B519" = "San Nicolás De Los Ranchos"
    "B521" = "San Salvador El Verde"
    "B528" = "Tlacotepec De Benito Juárez"
    "B539" = "Amealco De Bonfil"
    "B541" = "Cadereyta De Montes"
    "B544" = "Jalpan De Serra"
    "B545" = "Landa De Matamoros"
    "B546" = "Pinal De Amoles"
    "B548" = "San Juan Del Río"
    "B558" = "Ciudad Del Maíz"
    "B566" = "Santa María Del Río"
    "B570" = "Villa De Ramos"
    "B610" = "Soto La Marina"
    "B623" = "Alto Lucero De Gutiérrez Barrios"
    "B624" = "Amatlán De Los Reyes"
    "B640" = "Ignacio De La Llave"
    "B643" = "Juchique De Ferrer"
    "B646" = "Martínez De La Torre"
    "B649" = "Nanchital De Lázaro Cárdenas Del Río"
    "B654" = "Paso De Ovejas"
    "B655" = "Paso Del Macho"
    "B657" = "Poza Rica De Hidalgo"
    "B662" = "Sayula De Alemán"
    "B663" = "Soledad De Doblado"
    "B676" = "Vega De Alatorre"
    "B688" = "Concepción Del Oro"
    "B689" = "El Plateado De Joaquín Amaro"
    "B707" = "Tlaltenango De Sánchez Román"
}
foreach ($cellRef in $titleCaseFixes.Keys) {
    $ws.Range($cellRef).Value = $titleCaseFixes[$cellRef]
}

# --- Fix inconsistent capitalization typo ---
$ws.Range("B428").Value = "Montemorelos"

# --- Correct floating point rounding of a percentage value ---
$ws.Range("D457").Value = 0.009708737864077667

# --- Remove trailing footnote/metadata rows (716-720) ---
$ws.Range("A716:D720").ClearContents()
